# Auto-generated Excel COM-interop edit script
# Applies per-cell numeric updates (market price / profit recompute) to the
# 8 crafting-leve worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
#
# For each touched row we rewrite the changed currentAveragePrice* / Leve
# Price* / LeveProfit* cells (columns H-N). A handful of rows also gain or
# lose a trailing LeveProfitNQ/LeveProfitHQ cell (M/N) when a price feed
# flips between NQ-only and HQ-available data - those use ClearContents()
# to remove the now-absent cell instead of writing a value.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(42, 8).Value = 633.8  # H42: was 531.3333
$ws.Cells.Item(42, 9).Value = 723  # I42: was 547
$ws.Cells.Item(42, 11).Value = 2169  # K42: was 1641
$ws.Cells.Item(42, 13).Value = -1939  # M42: was -1411
$ws.Cells.Item(74, 8).Value = 13818.8  # H74: was 13819.3
$ws.Cells.Item(74, 9).Value = 16332  # I74: was 17285
$ws.Cells.Item(74, 10).Value = 6279.2  # J74: was 5732.6665
$ws.Cells.Item(74, 11).Value = 16332  # K74: was 17285
$ws.Cells.Item(74, 12).Value = 6279.2  # L74: was 5732.6665
$ws.Cells.Item(74, 13).Value = -15396  # M74: was -16349
$ws.Cells.Item(74, 14).Value = -8151.2  # N74: was -7604.6665
$ws.Cells.Item(77, 8).Value = 13818.8  # H77: was 13819.3
$ws.Cells.Item(77, 9).Value = 16332  # I77: was 17285
$ws.Cells.Item(77, 10).Value = 6279.2  # J77: was 5732.6665
$ws.Cells.Item(77, 11).Value = 81660  # K77: was 86425
$ws.Cells.Item(77, 12).Value = 31396  # L77: was 28663.3325
$ws.Cells.Item(77, 13).Value = -76980  # M77: was -81745
$ws.Cells.Item(77, 14).Value = -40756  # N77: was -38023.3325
$ws.Cells.Item(92, 8).Value = 431.30768  # H92: was 403.2143
$ws.Cells.Item(92, 9).Value = 431.30768  # I92: was 403.2143
$ws.Cells.Item(92, 11).Value = 431.30768  # K92: was 403.2143
$ws.Cells.Item(92, 13).Value = 816.69232  # M92: was 844.7857
$ws.Cells.Item(96, 8).Value = 174.6  # H96: was 161.25
$ws.Cells.Item(96, 9).Value = 95.111115  # I96: was 95
$ws.Cells.Item(96, 11).Value = 285.333345  # K96: was 285
$ws.Cells.Item(96, 13).Value = 1087.666655  # M96: was 1088
$ws.Cells.Item(111, 8).Value = 1155.2  # H111: was 1273.5
$ws.Cells.Item(111, 9).Value = 694  # I111: was 698
$ws.Cells.Item(111, 11).Value = 2082  # K111: was 2094
$ws.Cells.Item(111, 13).Value = 985  # M111: was 973
$ws.Cells.Item(125, 8).Value = 1500  # H125: was 1674
$ws.Cells.Item(125, 9).Value = 0  # I125: was 1698
$ws.Cells.Item(125, 10).Value = 1500  # J125: was 1650
$ws.Cells.Item(125, 11).Value = 0  # K125: was 15282
$ws.Cells.Item(125, 12).Value = 13500  # L125: was 14850
$ws.Cells.Item(125, 13).ClearContents()  # M125: was -12822
$ws.Cells.Item(125, 14).Value = -18420  # N125: was -19770
$ws.Cells.Item(137, 8).Value = 1519162.9  # H137: was 1519125
$ws.Cells.Item(137, 9).Value = 2003372.5  # I137: was 1926386.9
$ws.Cells.Item(137, 10).Value = 6007.875  # J137: was 6437.5713
$ws.Cells.Item(137, 11).Value = 6010117.5  # K137: was 5779160.699999999
$ws.Cells.Item(137, 12).Value = 18023.625  # L137: was 19312.7139
$ws.Cells.Item(137, 13).Value = -6007567.5  # M137: was -5776610.699999999
$ws.Cells.Item(137, 14).Value = -23123.625  # N137: was -24412.7139
$ws.Cells.Item(138, 8).Value = 2938.426  # H138: was 2941.6316
$ws.Cells.Item(138, 10).Value = 3645.3447  # J138: was 3584.7812
$ws.Cells.Item(138, 12).Value = 10936.0341  # L138: was 10754.3436
$ws.Cells.Item(138, 14).Value = -21216.0341  # N138: was -21034.3436

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 325  # H5: was 192
$ws.Cells.Item(5, 9).Value = 333.16666  # I5: was 165.66667
$ws.Cells.Item(5, 10).Value = 300.5  # J5: was 271
$ws.Cells.Item(5, 11).Value = 333.16666  # K5: was 165.66667
$ws.Cells.Item(5, 12).Value = 300.5  # L5: was 271
$ws.Cells.Item(5, 13).Value = -221.16666  # M5: was -53.66667000000001
$ws.Cells.Item(5, 14).Value = -524.5  # N5: was -495
$ws.Cells.Item(61, 8).Value = 3146.4211  # H61: was 2621.889
$ws.Cells.Item(61, 9).Value = 1945.75  # I61: was 1472.1818
$ws.Cells.Item(61, 10).Value = 3466.6  # J61: was 3412.3125
$ws.Cells.Item(61, 11).Value = 1945.75  # K61: was 1472.1818
$ws.Cells.Item(61, 12).Value = 3466.6  # L61: was 3412.3125
$ws.Cells.Item(61, 13).Value = -1733.75  # M61: was -1260.1818
$ws.Cells.Item(61, 14).Value = -3890.6  # N61: was -3836.3125
$ws.Cells.Item(63, 8).Value = 3157.375  # H63: was 3258.5715
$ws.Cells.Item(63, 9).Value = 2465.7144  # I63: was 2468.5
$ws.Cells.Item(63, 11).Value = 2465.7144  # K63: was 2468.5
$ws.Cells.Item(63, 13).Value = -1779.7144  # M63: was -1782.5
$ws.Cells.Item(66, 8).Value = 3157.375  # H66: was 3258.5715
$ws.Cells.Item(66, 9).Value = 2465.7144  # I66: was 2468.5
$ws.Cells.Item(66, 11).Value = 12328.572  # K66: was 12342.5
$ws.Cells.Item(66, 13).Value = -8896.572  # M66: was -8910.5
$ws.Cells.Item(74, 8).Value = 310819.88  # H74: was 279832.7
$ws.Cells.Item(74, 9).Value = 429212.16  # I74: was 372110.28
$ws.Cells.Item(74, 11).Value = 429212.16  # K74: was 372110.28
$ws.Cells.Item(74, 13).Value = -428338.16  # M74: was -371236.28
$ws.Cells.Item(77, 8).Value = 310819.88  # H77: was 279832.7
$ws.Cells.Item(77, 9).Value = 429212.16  # I77: was 372110.28
$ws.Cells.Item(77, 11).Value = 2146060.8  # K77: was 1860551.4
$ws.Cells.Item(77, 13).Value = -2141692.8  # M77: was -1856183.4
$ws.Cells.Item(97, 8).Value = 1238.6428  # H97: was 1210.4286
$ws.Cells.Item(97, 9).Value = 1137.8695  # I97: was 1124.2084
$ws.Cells.Item(97, 10).Value = 1702.2  # J97: was 1727.75
$ws.Cells.Item(97, 11).Value = 1137.8695  # K97: was 1124.2084
$ws.Cells.Item(97, 12).Value = 1702.2  # L97: was 1727.75
$ws.Cells.Item(97, 13).Value = -641.8695  # M97: was -628.2084
$ws.Cells.Item(97, 14).Value = -2694.2  # N97: was -2719.75
$ws.Cells.Item(102, 8).Value = 3199.2222  # H102: was 2744.3635
$ws.Cells.Item(102, 9).Value = 2942.875  # I102: was 2493.8
$ws.Cells.Item(102, 11).Value = 2942.875  # K102: was 2493.8
$ws.Cells.Item(102, 13).Value = -1320.875  # M102: was -871.8000000000002
$ws.Cells.Item(136, 8).Value = 3146.4211  # H136: was 2621.889
$ws.Cells.Item(136, 9).Value = 1945.75  # I136: was 1472.1818
$ws.Cells.Item(136, 10).Value = 3466.6  # J136: was 3412.3125
$ws.Cells.Item(136, 11).Value = 5837.25  # K136: was 4416.5454
$ws.Cells.Item(136, 12).Value = 10399.8  # L136: was 10236.9375
$ws.Cells.Item(136, 13).Value = -3287.25  # M136: was -1866.5454
$ws.Cells.Item(136, 14).Value = -15499.8  # N136: was -15336.9375

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 325  # H4: was 192
$ws.Cells.Item(4, 9).Value = 333.16666  # I4: was 165.66667
$ws.Cells.Item(4, 10).Value = 300.5  # J4: was 271
$ws.Cells.Item(4, 11).Value = 333.16666  # K4: was 165.66667
$ws.Cells.Item(4, 12).Value = 300.5  # L4: was 271
$ws.Cells.Item(4, 13).Value = -218.16666  # M4: was -50.66667000000001
$ws.Cells.Item(4, 14).Value = -530.5  # N4: was -501
$ws.Cells.Item(20, 8).Value = 18522708  # H20: was 22732256
$ws.Cells.Item(20, 9).Value = 29416762  # I20: was 38467840
$ws.Cells.Item(20, 10).Value = 2816  # J20: was 3078.889
$ws.Cells.Item(20, 11).Value = 29416762  # K20: was 38467840
$ws.Cells.Item(20, 12).Value = 2816  # L20: was 3078.889
$ws.Cells.Item(20, 13).Value = -29416515  # M20: was -38467593
$ws.Cells.Item(20, 14).Value = -3310  # N20: was -3572.889
$ws.Cells.Item(94, 8).Value = 80005736  # H94: was 74079470
$ws.Cells.Item(94, 9).Value = 100006300  # I94: was 90914930
$ws.Cells.Item(94, 11).Value = 100006300  # K94: was 90914930
$ws.Cells.Item(94, 13).Value = -100005849  # M94: was -90914479
$ws.Cells.Item(99, 8).Value = 94308.91  # H99: was 86649.586
$ws.Cells.Item(99, 9).Value = 127613  # I99: was 113678
$ws.Cells.Item(99, 10).Value = 5498  # J99: was 5564.3335
$ws.Cells.Item(99, 11).Value = 127613  # K99: was 113678
$ws.Cells.Item(99, 12).Value = 5498  # L99: was 5564.3335
$ws.Cells.Item(99, 13).Value = -126115  # M99: was -112180
$ws.Cells.Item(99, 14).Value = -8494  # N99: was -8560.333500000001
$ws.Cells.Item(105, 8).Value = 13001699  # H105: was 12382582
$ws.Cells.Item(105, 9).Value = 1001336.5  # I105: was 910328.6
$ws.Cells.Item(105, 11).Value = 1001336.5  # K105: was 910328.6
$ws.Cells.Item(105, 13).Value = -999589.5  # M105: was -908581.6
$ws.Cells.Item(134, 8).Value = 3199.7188  # H134: was 3151.0625
$ws.Cells.Item(134, 9).Value = 2539.5454  # I134: was 2470.087
$ws.Cells.Item(134, 10).Value = 4652.1  # J134: was 4891.3335
$ws.Cells.Item(134, 11).Value = 7618.6362  # K134: was 7410.261
$ws.Cells.Item(134, 12).Value = 13956.3  # L134: was 14674.0005
$ws.Cells.Item(134, 13).Value = -5083.6362  # M134: was -4875.261
$ws.Cells.Item(134, 14).Value = -19026.3  # N134: was -19744.0005

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 578.8  # H7: was 544.36365
$ws.Cells.Item(7, 9).Value = 448.25  # I7: was 420.66666
$ws.Cells.Item(7, 11).Value = 448.25  # K7: was 420.66666
$ws.Cells.Item(7, 13).Value = -335.25  # M7: was -307.66666
$ws.Cells.Item(31, 8).Value = 2782328.2  # H31: was 2721890.5
$ws.Cells.Item(31, 10).Value = 6255044  # J31: was 5957289.5
$ws.Cells.Item(31, 12).Value = 6255044  # L31: was 5957289.5
$ws.Cells.Item(31, 14).Value = -6255634  # N31: was -5957879.5
$ws.Cells.Item(34, 8).Value = 2782328.2  # H34: was 2721890.5
$ws.Cells.Item(34, 10).Value = 6255044  # J34: was 5957289.5
$ws.Cells.Item(34, 12).Value = 6255044  # L34: was 5957289.5
$ws.Cells.Item(34, 14).Value = -6255448  # N34: was -5957693.5
$ws.Cells.Item(64, 8).Value = 0  # H64: was 50000
$ws.Cells.Item(64, 9).Value = 0  # I64: was 50000
$ws.Cells.Item(64, 11).Value = 0  # K64: was 50000
$ws.Cells.Item(64, 13).ClearContents()  # M64: was -49752
$ws.Cells.Item(67, 8).Value = 0  # H67: was 50000
$ws.Cells.Item(67, 9).Value = 0  # I67: was 50000
$ws.Cells.Item(67, 11).Value = 0  # K67: was 50000
$ws.Cells.Item(67, 13).ClearContents()  # M67: was -49142
$ws.Cells.Item(107, 8).Value = 4167434  # H107: was 3846903.5
$ws.Cells.Item(107, 9).Value = 6250389  # I107: was 5555961
$ws.Cells.Item(107, 11).Value = 6250389  # K107: was 5555961
$ws.Cells.Item(107, 13).Value = -6248469  # M107: was -5554041

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 769.2857  # H107: was 723.125
$ws.Cells.Item(107, 10).Value = 834.25  # J107: was 747.4
$ws.Cells.Item(107, 12).Value = 2502.75  # L107: was 2242.2
$ws.Cells.Item(107, 14).Value = -6342.75  # N107: was -6082.2
$ws.Cells.Item(138, 8).Value = 7160.5264  # H138: was 6952.5
$ws.Cells.Item(138, 9).Value = 3256.25  # I138: was 3227.7778
$ws.Cells.Item(138, 11).Value = 9768.75  # K138: was 9683.3334
$ws.Cells.Item(138, 13).Value = -4628.75  # M138: was -4543.3334

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(15, 8).Value = 16698333  # H15: was 16683000
$ws.Cells.Item(15, 9).Value = 50000000  # I15: was 16683000
$ws.Cells.Item(15, 10).Value = 47499  # J15: was 0
$ws.Cells.Item(15, 11).Value = 50000000  # K15: was 16683000
$ws.Cells.Item(15, 12).Value = 47499  # L15: was 0
$ws.Cells.Item(15, 13).Value = -49999712  # M15: was -16682712
$ws.Cells.Item(15, 14).Value = -48075  # N15: was None
$ws.Cells.Item(81, 8).Value = 16698333  # H81: was 16683000
$ws.Cells.Item(81, 9).Value = 50000000  # I81: was 16683000
$ws.Cells.Item(81, 10).Value = 47499  # J81: was 0
$ws.Cells.Item(81, 11).Value = 50000000  # K81: was 16683000
$ws.Cells.Item(81, 12).Value = 47499  # L81: was 0
$ws.Cells.Item(81, 13).Value = -49999002  # M81: was -16682002
$ws.Cells.Item(81, 14).Value = -49495  # N81: was None
$ws.Cells.Item(84, 8).Value = 16698333  # H84: was 16683000
$ws.Cells.Item(84, 9).Value = 50000000  # I84: was 16683000
$ws.Cells.Item(84, 10).Value = 47499  # J84: was 0
$ws.Cells.Item(84, 11).Value = 150000000  # K84: was 50049000
$ws.Cells.Item(84, 12).Value = 142497  # L84: was 0
$ws.Cells.Item(84, 13).Value = -149995008  # M84: was -50044008
$ws.Cells.Item(84, 14).Value = -152481  # N84: was None
$ws.Cells.Item(102, 8).Value = 5466.6895  # H102: was 5532.241
$ws.Cells.Item(102, 9).Value = 1959.2084  # I102: was 2038.4166
$ws.Cells.Item(102, 11).Value = 1959.2084  # K102: was 2038.4166
$ws.Cells.Item(102, 13).Value = -337.2084  # M102: was -416.4166

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 6649.2856  # H7: was 6318.125
$ws.Cells.Item(7, 9).Value = 4710  # I7: was 4591.6665
$ws.Cells.Item(7, 11).Value = 4710  # K7: was 4591.6665
$ws.Cells.Item(7, 13).Value = -4598  # M7: was -4479.6665
$ws.Cells.Item(22, 8).Value = 400  # H22: was 1450
$ws.Cells.Item(22, 10).Value = 0  # J22: was 2500
$ws.Cells.Item(22, 12).Value = 0  # L22: was 2500
$ws.Cells.Item(22, 14).ClearContents()  # N22: was -3090
$ws.Cells.Item(27, 8).Value = 400  # H27: was 1450
$ws.Cells.Item(27, 10).Value = 0  # J27: was 2500
$ws.Cells.Item(27, 12).Value = 0  # L27: was 2500
$ws.Cells.Item(27, 14).ClearContents()  # N27: was -2714
$ws.Cells.Item(30, 8).Value = 5000  # H30: was 9999
$ws.Cells.Item(30, 9).Value = 5000  # I30: was 0
$ws.Cells.Item(30, 10).Value = 0  # J30: was 9999
$ws.Cells.Item(30, 11).Value = 5000  # K30: was 0
$ws.Cells.Item(30, 12).Value = 0  # L30: was 9999
$ws.Cells.Item(30, 13).Value = -4892  # M30: was None
$ws.Cells.Item(30, 14).ClearContents()  # N30: was -10215
$ws.Cells.Item(55, 8).Value = 754.7778  # H55: was 824.125
$ws.Cells.Item(55, 9).Value = 328.2  # I55: was 360.25
$ws.Cells.Item(55, 11).Value = 328.2  # K55: was 360.25
$ws.Cells.Item(55, 13).Value = -155.2  # M55: was -187.25
$ws.Cells.Item(93, 8).Value = 2363.2144  # H93: was 2079.0557
$ws.Cells.Item(93, 9).Value = 2335.75  # I93: was 1873.6666
$ws.Cells.Item(93, 10).Value = 2399.8333  # J93: was 2489.8333
$ws.Cells.Item(93, 11).Value = 2335.75  # K93: was 1873.6666
$ws.Cells.Item(93, 12).Value = 2399.8333  # L93: was 2489.8333
$ws.Cells.Item(93, 13).Value = -1087.75  # M93: was -625.6666
$ws.Cells.Item(93, 14).Value = -4895.8333  # N93: was -4985.8333
$ws.Cells.Item(100, 8).Value = 4854.4  # H100: was 3492.3572
$ws.Cells.Item(100, 9).Value = 3828.1428  # I100: was 2866.1667
$ws.Cells.Item(100, 10).Value = 7249  # J100: was 7249.5
$ws.Cells.Item(100, 11).Value = 3828.1428  # K100: was 2866.1667
$ws.Cells.Item(100, 12).Value = 7249  # L100: was 7249.5
$ws.Cells.Item(100, 13).Value = -3287.1428  # M100: was -2325.1667
$ws.Cells.Item(100, 14).Value = -8331  # N100: was -8331.5
$ws.Cells.Item(126, 8).Value = 6649.2856  # H126: was 6318.125
$ws.Cells.Item(126, 9).Value = 4710  # I126: was 4591.6665
$ws.Cells.Item(126, 11).Value = 14130  # K126: was 13774.9995
$ws.Cells.Item(126, 13).Value = -11660  # M126: was -11304.9995

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(51, 8).Value = 30070  # H51: was 0
$ws.Cells.Item(51, 9).Value = 30070  # I51: was 0
$ws.Cells.Item(51, 11).Value = 30070  # K51: was 0
$ws.Cells.Item(51, 13).Value = -29560  # M51: was None
$ws.Cells.Item(62, 8).Value = 6784.1816  # H62: was 6973.8
$ws.Cells.Item(62, 9).Value = 6038.6665  # I62: was 6268.8
$ws.Cells.Item(62, 11).Value = 6038.6665  # K62: was 6268.8
$ws.Cells.Item(62, 13).Value = -5414.6665  # M62: was -5644.8
$ws.Cells.Item(65, 8).Value = 6784.1816  # H65: was 6973.8
$ws.Cells.Item(65, 9).Value = 6038.6665  # I65: was 6268.8
$ws.Cells.Item(65, 11).Value = 30193.3325  # K65: was 31344
$ws.Cells.Item(65, 13).Value = -27073.3325  # M65: was -28224
$ws.Cells.Item(96, 8).Value = 4300  # H96: was 3533.6667
$ws.Cells.Item(96, 9).Value = 4300  # I96: was 3533.6667
$ws.Cells.Item(96, 11).Value = 4300  # K96: was 3533.6667
$ws.Cells.Item(96, 13).Value = -2927  # M96: was -2160.6667
$ws.Cells.Item(100, 8).Value = 62500572  # H100: was 50000510
$ws.Cells.Item(100, 9).Value = 549.6667  # I100: was 483.16666
$ws.Cells.Item(100, 10).Value = 142857740  # J100: was 125000560
$ws.Cells.Item(100, 11).Value = 1099.3334  # K100: was 966.33332
$ws.Cells.Item(100, 12).Value = 285715480  # L100: was 250001120
$ws.Cells.Item(100, 13).Value = -558.3334  # M100: was -425.33332
$ws.Cells.Item(100, 14).Value = -285716562  # N100: was -250002202
$ws.Cells.Item(136, 8).Value = 160619.17  # H136: was 166475.38
$ws.Cells.Item(136, 9).Value = 9682.328  # I136: was 9665.852999999999
$ws.Cells.Item(136, 10).Value = 560929.9399999999  # J136: was 644744.4
$ws.Cells.Item(136, 11).Value = 29046.984  # K136: was 28997.559
$ws.Cells.Item(136, 12).Value = 1682789.82  # L136: was 1934233.2
$ws.Cells.Item(136, 13).Value = -26496.984  # M136: was -26447.559
$ws.Cells.Item(136, 14).Value = -1687889.82  # N136: was -1939333.2
